# Auto-update data + news
# Update the "TDSP_pct" (Debt Service Burden) row (row 4) with refreshed
# values from the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 11.256338
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "Jul 2025"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = 11.09353892682927
$ws.Range("H4").Value = 0.1178030000000003
$ws.Range("I4").Value = 0.01057616643481394
